# Update init file references for scenario 21 (was scenario 20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value  = "A21"
$ws.Range("D6").Value  = "B21"
$ws.Range("D7").Value  = "C21"
$ws.Range("D8").Value  = "G21"
$ws.Range("D9").Value  = "H21"
$ws.Range("D10").Value = "I21"
$ws.Range("D11").Value = "J21"

# Update the active selection to match the saved workbook state
$ws.Range("D11").Select() | Out-Null
